# Auto-generated edit script: updates market-price / profit columns (H-N)
# on the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets to match the
# scheduled runner refresh described in the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 139.54546
$ws.Cells.Item(5, 9).Value = 93.57143000000001
$ws.Cells.Item(5, 10).Value = 220
$ws.Cells.Item(5, 11).Value = 93.57143000000001
$ws.Cells.Item(5, 12).Value = 220
$ws.Cells.Item(5, 13).Value = 21.42856999999999
$ws.Cells.Item(5, 14).Value = -450

$ws.Cells.Item(17, 8).Value = 370610.78
$ws.Cells.Item(17, 10).Value = 370610.78
$ws.Cells.Item(17, 12).Value = 1111832.34
$ws.Cells.Item(17, 14).Value = -1112168.34

$ws.Cells.Item(32, 8).Value = 1507.9
$ws.Cells.Item(32, 9).Value = 399.5
$ws.Cells.Item(32, 10).Value = 1785
$ws.Cells.Item(32, 11).Value = 399.5
$ws.Cells.Item(32, 12).Value = 1785
$ws.Cells.Item(32, 13).Value = -73.5
$ws.Cells.Item(32, 14).Value = -2437

$ws.Cells.Item(40, 8).Value = 2146.4443
$ws.Cells.Item(40, 9).Value = 2707.6
$ws.Cells.Item(40, 10).Value = 1445
$ws.Cells.Item(40, 11).Value = 2707.6
$ws.Cells.Item(40, 12).Value = 1445
$ws.Cells.Item(40, 13).Value = -2532.6
$ws.Cells.Item(40, 14).Value = -1795

$ws.Cells.Item(107, 8).Value = 1653
$ws.Cells.Item(107, 9).Value = 300
$ws.Cells.Item(107, 10).Value = 3006
$ws.Cells.Item(107, 11).Value = 300
$ws.Cells.Item(107, 12).Value = 3006
$ws.Cells.Item(107, 13).Value = 1620
$ws.Cells.Item(107, 14).Value = -6846

$ws.Cells.Item(137, 8).Value = 2779501.5
$ws.Cells.Item(137, 9).Value = 7693704.5
$ws.Cells.Item(137, 10).Value = 1908.3914
$ws.Cells.Item(137, 11).Value = 23081113.5
$ws.Cells.Item(137, 12).Value = 5725.174199999999
$ws.Cells.Item(137, 13).Value = -23078563.5
$ws.Cells.Item(137, 14).Value = -10825.1742

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1614.0588
$ws.Cells.Item(2, 9).Value = 1614.0588
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 1614.0588
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 13).Value = -1501.0588
$ws.Cells.Item(2, 14).ClearContents()

$ws.Cells.Item(34, 8).Value = 5000
$ws.Cells.Item(34, 9).Value = 5000
$ws.Cells.Item(34, 11).Value = 5000
$ws.Cells.Item(34, 13).Value = -4729

$ws.Cells.Item(61, 8).Value = 47715532
$ws.Cells.Item(61, 9).Value = 55612300
$ws.Cells.Item(61, 10).Value = 334933.34
$ws.Cells.Item(61, 11).Value = 55612300
$ws.Cells.Item(61, 12).Value = 334933.34
$ws.Cells.Item(61, 13).Value = -55612088
$ws.Cells.Item(61, 14).Value = -335357.34

$ws.Cells.Item(74, 8).Value = 10957823
$ws.Cells.Item(74, 9).Value = 16734811
$ws.Cells.Item(74, 10).Value = 125970
$ws.Cells.Item(74, 11).Value = 16734811
$ws.Cells.Item(74, 12).Value = 125970
$ws.Cells.Item(74, 13).Value = -16733937
$ws.Cells.Item(74, 14).Value = -127718

$ws.Cells.Item(77, 8).Value = 10957823
$ws.Cells.Item(77, 9).Value = 16734811
$ws.Cells.Item(77, 10).Value = 125970
$ws.Cells.Item(77, 11).Value = 83674055
$ws.Cells.Item(77, 12).Value = 629850
$ws.Cells.Item(77, 13).Value = -83669687
$ws.Cells.Item(77, 14).Value = -638586

$ws.Cells.Item(116, 8).Value = 1614.0588
$ws.Cells.Item(116, 9).Value = 1614.0588
$ws.Cells.Item(116, 10).Value = 0
$ws.Cells.Item(116, 11).Value = 1614.0588
$ws.Cells.Item(116, 12).Value = 0
$ws.Cells.Item(116, 13).Value = 679.9412
$ws.Cells.Item(116, 14).ClearContents()

$ws.Cells.Item(132, 8).Value = 80187.7
$ws.Cells.Item(132, 9).Value = 61368.766
$ws.Cells.Item(132, 10).Value = 112179.9
$ws.Cells.Item(132, 11).Value = 184106.298
$ws.Cells.Item(132, 12).Value = 336539.7
$ws.Cells.Item(132, 13).Value = -181576.298
$ws.Cells.Item(132, 14).Value = -341599.7

$ws.Cells.Item(136, 8).Value = 47715532
$ws.Cells.Item(136, 9).Value = 55612300
$ws.Cells.Item(136, 10).Value = 334933.34
$ws.Cells.Item(136, 11).Value = 166836900
$ws.Cells.Item(136, 12).Value = 1004800.02
$ws.Cells.Item(136, 13).Value = -166834350
$ws.Cells.Item(136, 14).Value = -1009900.02

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1614.0588
$ws.Cells.Item(3, 9).Value = 1614.0588
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 1614.0588
$ws.Cells.Item(3, 12).Value = 0
$ws.Cells.Item(3, 13).Value = -1500.0588
$ws.Cells.Item(3, 14).ClearContents()

$ws.Cells.Item(35, 8).Value = 10000
$ws.Cells.Item(35, 9).Value = 10000
$ws.Cells.Item(35, 11).Value = 10000
$ws.Cells.Item(35, 13).Value = -9690

$ws.Cells.Item(134, 8).Value = 1949.2609
$ws.Cells.Item(134, 9).Value = 2390.2307
$ws.Cells.Item(134, 10).Value = 1376
$ws.Cells.Item(134, 11).Value = 7170.6921
$ws.Cells.Item(134, 12).Value = 4128
$ws.Cells.Item(134, 13).Value = -4635.6921
$ws.Cells.Item(134, 14).Value = -9198

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 5701.25
$ws.Cells.Item(31, 9).Value = 2241.875
$ws.Cells.Item(31, 10).Value = 12620
$ws.Cells.Item(31, 11).Value = 2241.875
$ws.Cells.Item(31, 12).Value = 12620
$ws.Cells.Item(31, 13).Value = -1946.875
$ws.Cells.Item(31, 14).Value = -13210

$ws.Cells.Item(34, 8).Value = 5701.25
$ws.Cells.Item(34, 9).Value = 2241.875
$ws.Cells.Item(34, 10).Value = 12620
$ws.Cells.Item(34, 11).Value = 2241.875
$ws.Cells.Item(34, 12).Value = 12620
$ws.Cells.Item(34, 13).Value = -2039.875
$ws.Cells.Item(34, 14).Value = -13024

$ws.Cells.Item(58, 8).Value = 30304978
$ws.Cells.Item(58, 9).Value = 62501236
$ws.Cells.Item(58, 10).Value = 2615.647
$ws.Cells.Item(58, 11).Value = 62501236
$ws.Cells.Item(58, 12).Value = 2615.647
$ws.Cells.Item(58, 13).Value = -62501033
$ws.Cells.Item(58, 14).Value = -3021.647

$ws.Cells.Item(132, 8).Value = 42252.76
$ws.Cells.Item(132, 9).Value = 2583.8
$ws.Cells.Item(132, 10).Value = 101756.2
$ws.Cells.Item(132, 11).Value = 7751.400000000001
$ws.Cells.Item(132, 12).Value = 305268.6
$ws.Cells.Item(132, 13).Value = -5221.400000000001
$ws.Cells.Item(132, 14).Value = -310328.6

$ws.Cells.Item(134, 8).Value = 30840.459
$ws.Cells.Item(134, 9).Value = 1923.12
$ws.Cells.Item(134, 10).Value = 91084.914
$ws.Cells.Item(134, 11).Value = 5769.36
$ws.Cells.Item(134, 12).Value = 273254.742
$ws.Cells.Item(134, 13).Value = -3234.36
$ws.Cells.Item(134, 14).Value = -278324.742

$ws.Cells.Item(136, 8).Value = 30304978
$ws.Cells.Item(136, 9).Value = 62501236
$ws.Cells.Item(136, 10).Value = 2615.647
$ws.Cells.Item(136, 11).Value = 187503708
$ws.Cells.Item(136, 12).Value = 7846.941
$ws.Cells.Item(136, 13).Value = -187501158
$ws.Cells.Item(136, 14).Value = -12946.941

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(38, 8).Value = 145.82608
$ws.Cells.Item(38, 9).Value = 34.42857
$ws.Cells.Item(38, 10).Value = 194.5625
$ws.Cells.Item(38, 11).Value = 103.28571
$ws.Cells.Item(38, 12).Value = 583.6875
$ws.Cells.Item(38, 13).Value = 243.71429
$ws.Cells.Item(38, 14).Value = -1277.6875

$ws.Cells.Item(113, 8).Value = 561.76086
$ws.Cells.Item(113, 9).Value = 499.92
$ws.Cells.Item(113, 10).Value = 635.381
$ws.Cells.Item(113, 11).Value = 1499.76
$ws.Cells.Item(113, 12).Value = 1906.143
$ws.Cells.Item(113, 13).Value = 670.24
$ws.Cells.Item(113, 14).Value = -6246.143

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 107257.9
$ws.Cells.Item(132, 9).Value = 84503.164
$ws.Cells.Item(132, 10).Value = 146266
$ws.Cells.Item(132, 11).Value = 253509.492
$ws.Cells.Item(132, 12).Value = 438798
$ws.Cells.Item(132, 13).Value = -250979.492
$ws.Cells.Item(132, 14).Value = -443858

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 774.6667
$ws.Cells.Item(22, 9).Value = 550
$ws.Cells.Item(22, 10).Value = 887
$ws.Cells.Item(22, 11).Value = 550
$ws.Cells.Item(22, 12).Value = 887
$ws.Cells.Item(22, 13).Value = -255
$ws.Cells.Item(22, 14).Value = -1477

$ws.Cells.Item(27, 8).Value = 774.6667
$ws.Cells.Item(27, 9).Value = 550
$ws.Cells.Item(27, 10).Value = 887
$ws.Cells.Item(27, 11).Value = 550
$ws.Cells.Item(27, 12).Value = 887
$ws.Cells.Item(27, 13).Value = -443
$ws.Cells.Item(27, 14).Value = -1101

$ws.Cells.Item(46, 8).Value = 1422.1428
$ws.Cells.Item(46, 9).Value = 1633.3334
$ws.Cells.Item(46, 10).Value = 1364.5454
$ws.Cells.Item(46, 11).Value = 1633.3334
$ws.Cells.Item(46, 12).Value = 1364.5454
$ws.Cells.Item(46, 13).Value = -1445.3334
$ws.Cells.Item(46, 14).Value = -1740.5454

$ws.Cells.Item(48, 8).Value = 11500
$ws.Cells.Item(48, 9).Value = 11500
$ws.Cells.Item(48, 11).Value = 11500
$ws.Cells.Item(48, 13).Value = -10839

$ws.Cells.Item(68, 8).Value = 2582
$ws.Cells.Item(68, 9).Value = 1302
$ws.Cells.Item(68, 10).Value = 3222
$ws.Cells.Item(68, 11).Value = 1302
$ws.Cells.Item(68, 12).Value = 3222
$ws.Cells.Item(68, 13).Value = -553
$ws.Cells.Item(68, 14).Value = -4720

$ws.Cells.Item(71, 8).Value = 2582
$ws.Cells.Item(71, 9).Value = 1302
$ws.Cells.Item(71, 10).Value = 3222
$ws.Cells.Item(71, 11).Value = 6510
$ws.Cells.Item(71, 12).Value = 16110
$ws.Cells.Item(71, 13).Value = -2766
$ws.Cells.Item(71, 14).Value = -23598

$ws.Cells.Item(132, 8).Value = 38920.332
$ws.Cells.Item(132, 9).Value = 1193.3
$ws.Cells.Item(132, 10).Value = 146711.86
$ws.Cells.Item(132, 11).Value = 3579.9
$ws.Cells.Item(132, 12).Value = 440135.58
$ws.Cells.Item(132, 13).Value = -1049.9
$ws.Cells.Item(132, 14).Value = -445195.58

$ws.Cells.Item(136, 8).Value = 154800.77
$ws.Cells.Item(136, 9).Value = 167172.5
$ws.Cells.Item(136, 10).Value = 144196.42
$ws.Cells.Item(136, 11).Value = 501517.5
$ws.Cells.Item(136, 12).Value = 432589.26
$ws.Cells.Item(136, 13).Value = -498967.5
$ws.Cells.Item(136, 14).Value = -437689.26

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 119461.82
$ws.Cells.Item(132, 9).Value = 92281.17999999999
$ws.Cells.Item(132, 11).Value = 276843.54
$ws.Cells.Item(132, 13).Value = -274313.54

$ws.Cells.Item(136, 8).Value = 38691.395
$ws.Cells.Item(136, 9).Value = 22237.127
$ws.Cells.Item(136, 10).Value = 167583.17
$ws.Cells.Item(136, 11).Value = 66711.38099999999
$ws.Cells.Item(136, 12).Value = 502749.51
$ws.Cells.Item(136, 13).Value = -64161.38099999999
$ws.Cells.Item(136, 14).Value = -507849.51
